# The "conversion_factors_metal" sheet originally listed 9 metals
# (Mercury, Cadmium, Lead, Zinc, Copper, Nickel, Vanadium, Silver,
# Chromium) in rows 2-10. The edit keeps only 5 metals (Cadmium, Cuivre,
# Mercure, Plomb, Zinc), renames some labels to French, reorders them,
# and shrinks the used range down to A1:E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for Nickel, Vanadium, Silver and Chromium (rows 7-10).
# Delete from the bottom up so row indices of earlier rows stay valid.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()

# Re-populate rows 2-6 with the new ordering/labels, keeping each metal's
# original mussel/oyster/factor values attached to its new row position.
$ws.Range("A2").Value = "Cadmium"
$ws.Range("B2").Value = 1.124
$ws.Range("C2").Value = 1.5665
$ws.Range("D2").Value = 0.7175231407596554
$ws.Range("E2").Value = 1.393683274021352

$ws.Range("A3").Value = "Cuivre"
$ws.Range("B3").Value = 9.407499999999999
$ws.Range("C3").Value = 260.242
$ws.Range("D3").Value = 0.03614904588805803
$ws.Range("E3").Value = 27.6632474089822

$ws.Range("A4").Value = "Mercure"
$ws.Range("B4").Value = 0.1395
$ws.Range("C4").Value = 0.22
$ws.Range("D4").Value = 0.6340909090909091
$ws.Range("E4").Value = 1.577060931899641

$ws.Range("A5").Value = "Plomb"
$ws.Range("B5").Value = 1.24325
$ws.Range("C5").Value = 1.384
$ws.Range("D5").Value = 0.8983020231213873
$ws.Range("E5").Value = 1.11321134124271

$ws.Range("A6").Value = "Zinc"
$ws.Range("B6").Value = 144.0575
$ws.Range("C6").Value = 2880
$ws.Range("D6").Value = 0.05001996527777778
$ws.Range("E6").Value = 19.99201707651459
